$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 334: new case count changes from 71 to 70
$ws.Range("C334").Value = 70

# Row 336: new case count changes from 86 to 85
$ws.Range("C336").Value = 85

# Row 337: new case count changes from 57 to 91
$ws.Range("C337").Value = 91

# Row 338: new case count changes from 17 to 97; one extra extra-hospital death recorded
$ws.Range("C338").Value = 97
$ws.Range("M338").Value = 1

# Row 339: previously empty placeholder row, now filled in with actual data
$ws.Range("C339").Value = 12
$ws.Range("E339").Value = 15
$ws.Range("F339").Value = 11
$ws.Range("G339").Value = 105
$ws.Range("L339").Value = 0
$ws.Range("M339").Value = 0

$wb.Application.CalculateFullRebuild()
